# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# For Price (D) values that look like plain numbers we force the cell to
# text first (NumberFormat "@") so Excel doesn't silently coerce strings
# such as "205.94" or "0.485" into floating point numbers, then reset the
# style back to "Normal" so no stray number-format style is left applied
# to the cell (matches the original workbook, which stored these as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.900.72'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.547.54'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.485'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.31'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0856'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.97%  '
$ws.Range('D12').Value = '1.767.78'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').Value = '1.550.88'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.512'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '26.880.15'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '213.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').Value = '0.0₃0682'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.86%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E22').Value = '  -2.26%  '
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('E24').Value = '  -3.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.90'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0459'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.56%  '
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.08%  '
$ws.Range('D33').Value = '1.362.31'
$ws.Range('E33').Value = '  -3.04%  '
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.962'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.98%  '
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0165'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.520'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.807'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.989'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('E44').Value = '  +1.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.53'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('D47').Value = '1.681.86'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0508'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('D50').Value = '0.0₇0967'
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0950'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.30%  '
